# "Added Output photos + added missing slides"
#
# A new slide ("Drone Programming") is inserted right before the existing
# "Pluto X API (MSP Protocol)" slide (slide 17), pushing that slide and
# everything after it down by one position. Every other slide keeps its
# original content.
#
# Implementation: duplicate slide 17 (so the duplicate - an exact copy of
# "Pluto X API (MSP Protocol)" - lands immediately after it, at position
# 18), then overwrite the *original* slide 17's title/body text with the
# new "Drone Programming" content. The net effect is:
#   17 -> Drone Programming (new)
#   18 -> Pluto X API (MSP Protocol)  (was 17)
#   19 -> STAMP Support                (was 18)
#   20 -> Cloud Architecture           (was 19)
#   21 -> Gallery                      (was 20)
#   22 -> Future Scope                 (was 21)
#   23 -> Conclusion / Any questions   (was 22)

$p = $ppt.ActivePresentation

$original = $p.Slides.Item(17)

# Duplicate slide 17 - the copy is inserted right after it (position 18)
# and keeps the untouched "Pluto X API (MSP Protocol)" content.
$duplicate = $original.Duplicate()

# Replace the title and body of the original slide 17 with the new
# "Drone Programming" slide content.
$original.Shapes.Item(1).TextFrame.TextRange.Text = "Drone Programming"
$original.Shapes.Item(2).TextFrame.TextRange.Text = "Drone picture with small piece of code"
